$wb = $excel.ActiveWorkbook

# The workbook was edited in the "Lookup Functions" sheet: a VLOOKUP array
# formula was entered in E2 and filled down to E205, looking up each state's
# postal abbreviation from the "State Abbreviations" sheet. Filling it in via
# the UI also makes "Lookup Functions" the active/selected sheet (previously
# "Stats Functions" was the selected tab).

$ws = $wb.Worksheets.Item("Lookup Functions")
$ws.Activate()

$ws.Range("E2:E205").FormulaArray = "=VLOOKUP(A2:A205,'State Abbreviations'!A2:B52,2)"

[void]$ws.Range("E3").Select()
